$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.307.42"
$ws.Range("E2").Value = "  +1.35%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.486.22"
$ws.Range("E3").Value = "  +3.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.64%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'578.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'147.02"
$ws.Range("D6").Style = "Normal"

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.486.19"
$ws.Range("E9").Value = "  +2.17%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.71%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.59%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +0.60%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'28.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.29%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.939.95"
$ws.Range("E16").Value = "  +1.95%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.239.59"
$ws.Range("E17").Value = "  +1.02%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.488.20"
$ws.Range("E18").Value = "  +3.00%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'8.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.63%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +1.23%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'330.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "

# Row 22 - SuiNetwork
$ws.Range("D22").Value = "'2.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.11%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.04%  "

# Row 25 - Litecoin
$ws.Range("E25").Value = "  +1.20%  "

# Row 26 - Aptos
$ws.Range("D26").Value = "'9.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.20%  "

# Row 27 - Bittensor
$ws.Range("D27").Value = "'667.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.49%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +2.12%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "2.608.04"
$ws.Range("E29").Value = "  +1.83%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +311.05%  "

# Row 31 - Fetch.AI
$ws.Range("D31").Value = "'1.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.68%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'8.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +0.32%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -2.88%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "'1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.66%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.28%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +1.02%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "'5.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("D39").Value = "'0.373"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

# Row 40 - EthereumClassic
$ws.Range("E40").Value = "  +0.96%  "

# Row 41 - Monero
$ws.Range("D41").Value = "'150.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.94%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "'2.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +1.30%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0313"
$ws.Range("E45").Value = "  -15.31%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'156.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.42%  "

# Row 47 - WhiteBITCoin
$ws.Range("D47").Value = "'15.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.47%  "

# Row 48 - Filecoin
$ws.Range("D48").Value = "'3.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'20.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "'0.609"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.91%  "

# Row 51 - Hedera
$ws.Range("D51").Value = "'0.0516"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
